# Harness Lab Presentation — replace the small "summary tables" next to
# screenshots with full-bleed screenshots (and drop the now-redundant
# caption text boxes on the Pipeline Execution slide).
#
# EMU -> point conversions below were chosen so that, after the COM
# layer's internal Single(32-bit float) round trip, the saved EMU in the
# OOXML lands back on the exact target value (PowerPoint COM stores
# Left/Top/Width/Height in points as a 32-bit float; 1 pt = 12700 EMU).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 6 — "Connectors — Configure Once, Use Everywhere"
# Remove the connector-status table and blow the screenshot up to fill
# the content area.
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(2).Delete()                 # "Table 2"

$s6pic = $s6.Shapes.Item(2)                 # was "Picture 3"
$s6pic.Name = "Picture 2"
$s6pic.Left = 36.0
$s6pic.Top = 79.2000008
$s6pic.Width = 885.6000062
$s6pic.Height = 422.5846456692913

# ---------------------------------------------------------------------
# Slide 7 — "CI/CD Pipeline — Visual + YAML"
# Remove the stage-by-stage bullet text box and blow the pipeline
# screenshot up to fill the content area.
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(2).Delete()                 # "TextBox 2"

$s7pic = $s7.Shapes.Item(2)                 # was "Picture 3"
$s7pic.Name = "Picture 2"
$s7pic.Left = 36.0
$s7pic.Top = 79.2000008
$s7pic.Width = 885.6000062
$s7pic.Height = 481.1067716535433

# ---------------------------------------------------------------------
# Slide 9 — "Pipeline Execution — Success"
# Drop the CI/CD results caption text boxes underneath the screenshot.
# (delete highest index first so the remaining indices don't shift)
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(4).Delete()                 # "TextBox 4" (CD Results)
$s9.Shapes.Item(3).Delete()                 # "TextBox 3" (CI Results)

# ---------------------------------------------------------------------
# Slide 10 — "Live Application — Deployed & Running"
# Remove the endpoint/response table, blow the screenshot up to fill
# the content area, and drop the live-URL caption further down with a
# larger font.
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(2).Delete()                # "Table 2"

$s10pic = $s10.Shapes.Item(2)               # was "Picture 3"
$s10pic.Name = "Picture 2"
$s10pic.Left = 36.0
$s10pic.Top = 93.6000023
$s10pic.Width = 885.6000062
$s10pic.Height = 298.4086762173228

$s10tb = $s10.Shapes.Item(3)                # was "TextBox 4"
$s10tb.Name = "TextBox 3"
$s10tb.Left = 36.0
$s10tb.Top = 396.0
$s10tb.Width = 864.0
$s10tb.Height = 72.0
$s10tb.TextFrame.TextRange.Font.Size = 16
